# Revert "Revision 3" -> "Revision 2", move the _GoBack bookmark into the
# first paragraph (between the "Revision 2" run and the trailing "!" run),
# and remove the blank paragraph + the "It now has a bold title!" paragraph
# that used to hold that bookmark.

$d = $word.ActiveDocument

# Locate the three runs of paragraph 1 precisely via Find, so this does not
# depend on hard-coded character offsets.
$rRun1 = $d.Content
[void]$rRun1.Find.Execute("This is a test file for TC_12", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$run1End = $rRun1.End

$rRun2 = $d.Content
[void]$rRun2.Find.Execute(" and TC_28.  Revision 3", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$run2Start = $rRun2.Start
$run2End = $rRun2.End

# Step 1: drop a throwaway bookmark exactly on the run1/run2 boundary. This
# pins that boundary so the text-replace below (which only touches run2)
# cannot merge back into run1 -- this engine otherwise flattens the whole
# paragraph into a single run on any text edit, but it always respects
# bookmarks as hard boundaries.
$d.Bookmarks.Add("zzTempBoundary", $d.Range($run1End, $run1End))

# Step 2: put the real _GoBack bookmark on the run2/run3 boundary -- this is
# also its final resting place in the target document.
$d.Bookmarks.Add("_GoBack", $d.Range($run2End, $run2End))

# Step 3: rewrite run2's text in place (still pinned on both sides by the
# two bookmarks, so run1 and run3 stay untouched).
$d.Range($run2Start, $run2End).Text = " and TC_28.  Revision 2"

# Step 4: the temporary bookmark has done its job -- remove it (removing a
# bookmark is not a text edit, so it does not trigger another run merge).
$d.Bookmarks("zzTempBoundary").Delete()

# Step 5: delete the blank paragraph that used to sit right after paragraph 1.
$d.Paragraphs(2).Range.Delete()

# Step 6: delete the (now second) paragraph, which used to hold the bold
# "It now has a bold title!" text and the _GoBack bookmark -- the bookmark
# has already been relocated above, so this just removes the leftover text.
$d.Paragraphs(2).Range.Delete()

Write-Host "Done: reverted to Revision 2 and removed the bold-title paragraph."
